$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "IRS f1040 C"
